$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text in E8 (was "Good Morning", now "GIT UPDATE")
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new selection/active cell shown in the saved sheet view
$ws.Range("E8").Select()
